$wb = $excel.ActiveWorkbook

# Rename the existing "strategy_id-5008" sheet to "strategy_id-5007"
$src = $wb.Worksheets.Item("strategy_id-5008")
$src.Name = "strategy_id-5007"

# Duplicate the sheet (copies all formatting/values) and place the copy
# immediately after the renamed sheet, then rename the copy to
# "strategy_id-5009" per the new template list.
$src.Copy($null, $src)

$new = $wb.Worksheets.Item($src.Index + 1)
$new.Name = "strategy_id-5009"
